# Edit: "dat field issues with mozilla" -> "date field issues with mozilla"
# (split into separate runs, proofErr moves from "dat" to "mozilla"), then add
# a blank paragraph and a new paragraph "added task create and mofification
# code, improved document form further with more functionality" (with the
# _GoBack bookmark moved to the end of that new paragraph).

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1) Replace the final paragraph's content ("dat field issues with mozilla")
#    with the reworded/re-run version ("date field issues with mozilla"),
#    dropping the old bookmark (it will be re-added to the new last paragraph
#    below) and moving the proofErr spell-check wrapper from "dat" onto
#    "mozilla".
$target = $d.Paragraphs.Last.Range
$null = $target.InsertXML("<w:p xmlns:w=""$wNs""><w:r><w:t>dat</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space=""preserve""> field issues with </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>mozilla</w:t></w:r><w:proofErr w:type=""spellEnd""/></w:p>")

# 2) Insert a new blank paragraph right after it.
$target = $d.Paragraphs.Last.Range
$null = $target.InsertParagraphAfter()
$blank = $d.Paragraphs.Last.Range
$null = $blank.InsertXML("<w:p xmlns:w=""$wNs""/>")

# 3) Insert the new trailing paragraph (with the _GoBack bookmark moved here).
$blank = $d.Paragraphs.Last.Range
$null = $blank.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last.Range
$null = $newPara.InsertXML("<w:p xmlns:w=""$wNs""><w:r><w:t xml:space=""preserve"">added task create and </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>mofification</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> code, improved document form further with more functionality</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>")
